# VIC_COVID19_DAILY_DATA.xlsx — append daily rows 101-108 (dates 44904-44911)
# and restyle the "New Cases"/"PCR cases"/"PCR Tests" columns with a
# thousands-separator number format, then freeze the header row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) New daily data rows (A:H), 2022-12-09 .. 2022-12-16
# ---------------------------------------------------------------------------
$data = @(
    @(44904, 3053, 2164, 705, 30, 21, 6899),
    @(44905, 2686, 1884, 719, 31, 0, 6267),
    @(44906, 2873, 2314, 687, 23, 0, 4534),
    @(44907, 4376, 3542, 680, 26, 20, 6642),
    @(44908, 3882, 2807, 660, 23, 10, 8654),
    @(44909, 4209, 3270, 639, 25, 20, 7445),
    @(44910, 3561, 2607, 638, 28, 17, 7294),
    @(44911, 3371, 2431, 644, 24, 19, 7638)
)

$row = 101
foreach ($r in $data) {
    $ws.Cells.Item($row, 1).Value = $r[0]   # A: Date reported
    $ws.Cells.Item($row, 2).Value = $r[1]   # B: New Cases
    # C (PCR cases) filled in below via the shared formula B-D
    $ws.Cells.Item($row, 4).Value = $r[2]   # D: RAT cases
    $ws.Cells.Item($row, 5).Value = $r[3]   # E: In Hospital
    $ws.Cells.Item($row, 6).Value = $r[4]   # F: In ICU
    $ws.Cells.Item($row, 7).Value = $r[5]   # G: Deaths
    $ws.Cells.Item($row, 8).Value = $r[6]   # H: PCR Tests
    $row++
}

# Column C: "PCR cases" = New Cases - RAT cases, filled as one shared formula
$ws.Range("C101:C108").Formula = "=B101-D101"

# ---------------------------------------------------------------------------
# 2) Number formatting
# ---------------------------------------------------------------------------

# A101:A108 — same date format as the rest of column A (copy format only,
# so the freshly entered values/formula are left untouched).
$ws.Range("A100").Copy() | Out-Null
$ws.Range("A101:A108").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = $false

# B100/C100/H100 (pre-existing row) plus the whole new block B101:H108 get a
# thousands-separator number format (#,##0 = numFmtId 3).
$ws.Range("B100").NumberFormat = "#,##0"
$ws.Range("C100").NumberFormat = "#,##0"
$ws.Range("H100").NumberFormat = "#,##0"
$ws.Range("B101:B108").NumberFormat = "#,##0"
$ws.Range("C101:C108").NumberFormat = "#,##0"
$ws.Range("H101:H108").NumberFormat = "#,##0"

# ---------------------------------------------------------------------------
# 3) View: freeze the header row, scroll near the bottom, select the last cell
# ---------------------------------------------------------------------------
$ws.Range("A2").Select()
$win = $excel.ActiveWindow
$win.FreezePanes = $true

$ws.Range("A108").Select()
